# Update the "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets. These two sheets carry duplicate data, so the same cells need
# to be updated in both places.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7692
    9  = 5905
    14 = 1315
    16 = 447
    17 = 104
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
